# Community Lab - Intro to Composable (Common)
# Commit: "fix spelling error in slide"
#
# 1) Fix the typo "afe" -> "are" in the "@staticmethod provides means..."
#    sentence on slide 9 (without disturbing the surrounding run/paragraph
#    structure - the run is replaced in full so PowerPoint keeps a single
#    <a:r> instead of splitting it into three runs).
# 2) The underlying OOXML also shows the deck's two theme parts trading
#    their <a:clrScheme> colors (the color scheme actually used by the
#    slide master moves from the "Simple Light" palette to the "Default"
#    palette). Re-apply that palette through the Theme Color Scheme the
#    object model exposes for the active design.

$p = $ppt.ActivePresentation

# --- 1. Fix the typo -------------------------------------------------
$oldSentence = "@staticmethod provides means to tear off buildable micro components that afe configured by metaparameters (factory design pattern)."
$newSentence = "@staticmethod provides means to tear off buildable micro components that are configured by metaparameters (factory design pattern)."

$s  = $p.Slides.Item(9)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf($oldSentence)

if ($idx -lt 0) {
    # Fall back to a full scan in case shape/slide layout ever shifts.
    for ($si = 1; $si -le $p.Slides.Count -and $idx -lt 0; $si++) {
        $slide = $p.Slides.Item($si)
        for ($shi = 1; $shi -le $slide.Shapes.Count -and $idx -lt 0; $shi++) {
            $cand = $slide.Shapes.Item($shi)
            if ($cand.HasTextFrame) {
                $candTr = $cand.TextFrame.TextRange
                $candFull = $candTr.Text
                $candIdx = $candFull.IndexOf($oldSentence)
                if ($candIdx -ge 0) {
                    $tr = $candTr
                    $full = $candFull
                    $idx = $candIdx
                }
            }
        }
    }
}

if ($idx -ge 0) {
    $run = $tr.Characters($idx + 1, $oldSentence.Length)
    $run.Text = $newSentence
}

# --- 2. Swap the theme's colour scheme (Simple Light -> Default) -----
$tcs = $s.ThemeColorScheme

# Target palette = the "Default" clrScheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) expressed as COLORREF (0x00BBGGRR) ints.
$defaultPalette = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    5800213,   # dk2      158158
    15987699,  # lt2      F3F3F3
    13077765,  # accent1  058DC7
    3322960,   # accent2  50B432
    1791725,   # accent3  ED561B
    61421,     # accent4  EDEF00
    15059748,  # accent5  24CBE5
    7529828,   # accent6  64E572
    13369378,  # hlink    2200CC
    9116245    # folHlink 551A8B
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $defaultPalette[$i - 1]
}
